# "integrate files from next team"
#
# Applies the substantive content edits captured in the target diff:
#   - About sheet: drop the Realmonte et al. citation block and the
#     "model end date" notes block, replace with new, shorter notes
#     about South Korea not using DAC and scaling the global data.
#   - Data sheet: U.S. GDP figure revised (19.39 -> 1.624 trillion USD)
#     and a new "GDP from WorldBank" source note alongside it.
#
# Everything downstream (Data!B74, rows 78/79/83/84, and the TREND()
# cells on DACD-potential) is a plain formula cascade, so we just let
# the workbook recalc after the writes instead of hand-computing them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# About sheet
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Drop the hyperlinks that point at the Realmonte et al. source material
# before clearing the cells that hosted them.
$about.Range("B9").Hyperlinks.Delete()
$about.Range("B11").Hyperlinks.Delete()

# Citation block (Realmonte et al., the two source URLs, "Table 2", and
# the supplementary-materials caption) is replaced by a single "N/A" —
# no longer a cited source — and the rest of the block is cleared out.
$about.Range("B6").Value = "N/A"
$about.Range("B7").ClearContents()
$about.Range("B8").ClearContents()
$about.Range("B9").ClearContents()
$about.Range("B10").ClearContents()
$about.Range("B11").ClearContents()
$about.Range("B12").ClearContents()

# Notes block: the old "model end date" caveat is replaced with a note
# about South Korea not using DAC, and the long explanation of why only
# the hydroxide-sorbent technology is modeled is replaced with a note
# about scaling down the global data using the method used for the US.
$about.Range("A15").Value = "DAC is not included in any energy pathways in South Korea"
$about.Range("A16").Value = "Therefore we will scale down the global data following the method used for the US"
$about.Range("A17").ClearContents()
$about.Range("A18").ClearContents()
$about.Range("A20").ClearContents()
$about.Range("A21").ClearContents()
$about.Range("A22").ClearContents()
$about.Range("A23").ClearContents()

# Amortized CapEx and OM Cost Notes block is removed entirely.
$about.Range("A25").ClearContents()
$about.Range("A26").ClearContents()
$about.Range("A27").ClearContents()

# ---------------------------------------------------------------------
# Data sheet
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

# U.S. GDP (trillion USD, 2017) revised down from 19.39 to 1.624, and a
# source note added alongside the World GDP figure's existing source.
$data.Range("B72").Value = 1.624
$data.Range("E72").Value = "GDP from WorldBank"

$wb.Application.CalculateFull()
